$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All test cases are now enabled (EXECUTION STATE column set to "Y")
$ws.Range("A2").Value = "Y"
$ws.Range("A3").Value = "Y"
$ws.Range("A4").Value = "Y"
$ws.Range("A5").Value = "Y"
$ws.Range("A6").Value = "Y"
$ws.Range("A7").Value = "Y"
$ws.Range("A8").Value = "Y"
$ws.Range("A9").Value = "Y"
$ws.Range("A10").Value = "Y"
$ws.Range("A11").Value = "Y"

# test_case_10's environment value changed from 1 to 2
$ws.Range("D11").Value = 2

# Update the last active selection on the sheet
$ws.Range("E15").Select()
